# BadEventsV4ToParse.xlsx — "making the deifferent events text show up"
#
# 1. BAD!C11:I11 were empty placeholder cells for a row that had no text -
#    fill them with the literal "blank" so the row reads consistently.
# 2. BAD!B14 held the generic "BLANK" shared string (index 86, shared with
#    C14:I14) - give it its own real event text.
# 3. BAD!E15 and BAD!H15 were empty - fill them with "blank" too.
# 4. View/selection bookkeeping: the BAD sheet becomes the active/selected
#    tab (was GOOD), with a new active cell, and the GOOD sheet's own
#    selection moves off of B2 onto D2.

$wb = $excel.ActiveWorkbook

$wsBad = $wb.Worksheets.Item("BAD")
$wsGood = $wb.Worksheets.Item("GOOD")

# --- data edits -----------------------------------------------------------

# Row 11 was entirely blank (except col A); give the rest of the row
# explicit "blank" text so it shows up instead of looking empty/broken.
$wsBad.Range("C11").Value = "blank"
$wsBad.Range("D11").Value = "blank"
$wsBad.Range("E11").Value = "blank"
$wsBad.Range("F11").Value = "blank"
$wsBad.Range("G11").Value = "blank"
$wsBad.Range("H11").Value = "blank"
$wsBad.Range("I11").Value = "blank"

# B14 used to just repeat the generic "BLANK" placeholder - now it gets the
# real "different" event text that should actually show up.
$wsBad.Range("B14").Value = "You saw him trying to command your followers in giving him things. "

# E15 / H15 were also blank cells - fill with "blank" like row 11.
$wsBad.Range("E15").Value = "blank"
$wsBad.Range("H15").Value = "blank"

# --- view / selection edits -----------------------------------------------

# Move the selected tab from GOOD to BAD, scroll/select B14 on BAD.
[void]$wsBad.Activate()
[void]$excel.Goto($wsBad.Range("A7"), $false)
[void]$wsBad.Range("B14").Select()

# GOOD sheet's own remembered selection moves from B2 to D2.
[void]$wsGood.Range("D2").Select()

# Leave BAD as the active sheet/tab (matches tabSelected moving off GOOD).
[void]$wsBad.Activate()
[void]$wsBad.Range("B14").Select()
